$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "None of the file or data operations ... is performed ..." paragraph:
#    remove the proofErr gramStart/gramEnd wrapper around "is" by replacing
#    the paragraph's text with a single run.
# ---------------------------------------------------------------------------
$pNone = $d.Paragraphs.Item(3)
$rNone = $pNone.Range
$rNoneBody = $d.Range($rNone.Start, $rNone.End - 1)
$rNoneBody.Delete()
$insNone = $d.Range($rNone.Start, $rNone.Start)
$insNone.InsertAfter("None of the file or data operations is performed in the user interface, which is how novice developers write code, while a seasoned developer writes file and data operations in classes.")

# ---------------------------------------------------------------------------
# 2. "Note: EF Core rather than Dapper would be overkill." paragraph: remove
#    the proofErr gramStart/gramEnd wrapper around "would be" while keeping
#    the three runs distinct.
# ---------------------------------------------------------------------------
$pNote = $d.Paragraphs.Item(6)
$rNote = $pNote.Range
$rNoteBody = $d.Range($rNote.Start, $rNote.End - 1)
$rNoteBody.Delete()

$insNote1 = $d.Range($rNote.Start, $rNote.Start)
$insNote1.InsertAfter("Note: EF Core rather than Dapper ")

$start2 = $insNote1.End
$midHolder = $d.Range($start2, $start2)
$midHolder.InsertAfter("would be")
$midRange = $d.Range($start2, $start2 + 8)
# Toggle a character property so this run stays split from its neighbours
# (adjacent runs with identical formatting are coalesced on save).
$midRange.Bold = 1
$midRange.Bold = 0

$start3 = $midRange.End
$insNote3 = $d.Range($start3, $start3)
$insNote3.InsertAfter(" overkill.")

# ---------------------------------------------------------------------------
# 3. Insert two new bullet paragraphs between the "None of..." paragraph and
#    the "Step 1" heading: "Talk about " (top level) and "The preview
#    features configuration in the project file" (second level).
# ---------------------------------------------------------------------------
$pNone = $d.Paragraphs.Item(3)
$newPara1 = $pNone.Range.InsertParagraphAfter()
$pTalk = $d.Paragraphs.Item(4)
$pTalk.Range.Text = "Talk about "

$newPara2 = $pTalk.Range.InsertParagraphAfter()
$pPreview = $d.Paragraphs.Item(5)
$pPreview.Range.Text = "The preview features configuration in the project file"
$pPreview.Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------------
# 4. The second list level (ilvl=1) of the numbering definition used above is
#    now actually used in the document, so clear its "tentative" flag.
# ---------------------------------------------------------------------------
$pAnyListItem = $d.Paragraphs.Item(3)
$listTemplate = $pAnyListItem.Range.ListFormat.ListTemplate
$secondLevel = $listTemplate.ListLevels.Item(2)
$secondLevel.NumberFormat = $secondLevel.NumberFormat

Write-Output "done"
